# Update weekly fruit/vegetable (Pera asiática) price records
# Applies the diff mapping rows to their shuffled date/price data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = (Get-Date -Year 2023 -Month 2 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 250000
$ws.Range("O2").Value = 250000
$ws.Range("P2").Value = 250000
$ws.Range("Q2").Value = "`$/bins (400 kilos)"
$ws.Range("S2").Value = 625
$ws.Range("T2").Value = 400

# Row 3
$ws.Range("D3").Value = (Get-Date -Year 2023 -Month 2 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("S3").Value = 833

# Row 5
$ws.Range("D5").Value = (Get-Date -Year 2021 -Month 8 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M5").Value = 210
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("Q5").Value = "`$/bandeja 8 kilos"
$ws.Range("S5").Value = 1000
$ws.Range("T5").Value = 8

# Row 6
$ws.Range("D6").Value = (Get-Date -Year 2023 -Month 6 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 110
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("Q6").Value = "`$/caja 18 kilos granel"
$ws.Range("S6").Value = 889
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = (Get-Date -Year 2021 -Month 1 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 24000
$ws.Range("P7").Value = 24000
$ws.Range("S7").Value = 1600

# Row 8
$ws.Range("D8").Value = (Get-Date -Year 2023 -Month 6 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M8").Value = 55

# Row 9
$ws.Range("D9").Value = (Get-Date -Year 2022 -Month 2 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 28000
$ws.Range("O9").Value = 28000
$ws.Range("P9").Value = 28000
$ws.Range("Q9").Value = "`$/caja 18 kilos granel"
$ws.Range("S9").Value = 1556
$ws.Range("T9").Value = 18

# Row 11
$ws.Range("D11").Value = (Get-Date -Year 2021 -Month 3 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L11").Value = "Calibre 100"
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("Q11").Value = "`$/caja 18 kilos embalada"
$ws.Range("S11").Value = 1111
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = (Get-Date -Year 2023 -Month 6 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 16000
$ws.Range("Q12").Value = "`$/caja 18 kilos granel"
$ws.Range("S12").Value = 889

# Row 13
$ws.Range("D13").Value = (Get-Date -Year 2023 -Month 6 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M13").Value = 170
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 15471
$ws.Range("Q13").Value = "`$/caja 18 kilos granel"
$ws.Range("S13").Value = 860
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = (Get-Date -Year 2021 -Month 11 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M14").Value = 15
$ws.Range("N14").Value = 22000
$ws.Range("O14").Value = 22000
$ws.Range("P14").Value = 22000
$ws.Range("Q14").Value = "`$/caja 15 kilos granel"
$ws.Range("S14").Value = 1467
$ws.Range("T14").Value = 15

# Row 15
$ws.Range("D15").Value = (Get-Date -Year 2023 -Month 6 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M15").Value = 220
$ws.Range("N15").Value = 16000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 16000
$ws.Range("S15").Value = 889

# Row 16
$ws.Range("D16").Value = (Get-Date -Year 2021 -Month 7 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 500
$ws.Range("N16").Value = 7000
$ws.Range("O16").Value = 7000
$ws.Range("P16").Value = 7000
$ws.Range("Q16").Value = "`$/bandeja 8 kilos"
$ws.Range("S16").Value = 875
$ws.Range("T16").Value = 8

# Row 17
$ws.Range("D17").Value = (Get-Date -Year 2021 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M17").Value = 55
$ws.Range("N17").Value = 18000
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 18000
$ws.Range("Q17").Value = "`$/caja 18 kilos granel"
$ws.Range("T17").Value = 18

# Row 18
$ws.Range("D18").Value = (Get-Date -Year 2023 -Month 6 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 16000
$ws.Range("P18").Value = 16000
$ws.Range("S18").Value = 889

# Row 19
$ws.Range("D19").Value = (Get-Date -Year 2021 -Month 8 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L19").Value = "Especial"
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 8000
$ws.Range("O19").Value = 8000
$ws.Range("P19").Value = 8000
$ws.Range("Q19").Value = "`$/caja 15 kilos granel"
$ws.Range("S19").Value = 533
$ws.Range("T19").Value = 15

# Row 20
$ws.Range("D20").Value = (Get-Date -Year 2021 -Month 8 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N20").Value = 7000
$ws.Range("O20").Value = 7000
$ws.Range("P20").Value = 7000
$ws.Range("Q20").Value = "`$/caja 15 kilos granel"
$ws.Range("S20").Value = 467
$ws.Range("T20").Value = 15

